$d = $word.ActiveDocument

# Remove the first four paragraphs at the top of the form:
#   "Name: "
#   (empty)
#   "Poker Skill Level:"
#   (empty)
# Deleting the range from the start of the document through the end of
# the 4th paragraph (including its paragraph mark) leaves "Task # 1:"
# as the new first paragraph, matching the rest of the document.
$endOfFourthParagraph = $d.Paragraphs.Item(4).Range.End
$deleteRange = $d.Range(0, $endOfFourthParagraph)
$deleteRange.Delete()
